$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with the latest scraped figures. Price values that look numeric (e.g.
# "535.63") are written with a leading apostrophe so Excel keeps them as
# plain text, matching the rest of the (text-formatted) Price column
# instead of being auto-coerced into the Number type.
$ws.Range("D2").Value = '71.021.63'
$ws.Range("E2").Value = '  -1.86%  '
$ws.Range("D3").Value = '3.945.68'
$ws.Range("E3").Value = '  -2.55%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''535.63'
$ws.Range("E5").Value = '  +2.64%  '
$ws.Range("D6").Value = '''147.90'
$ws.Range("E6").Value = '  -0.44%  '
$ws.Range("D7").Value = '3.940.05'
$ws.Range("E7").Value = '  -2.53%  '
$ws.Range("E8").Value = '  -4.49%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").Value = '''0.739'
$ws.Range("E10").Value = '  -5.27%  '
$ws.Range("D11").Value = '''0.166'
$ws.Range("E11").Value = '  -7.02%  '
$ws.Range("D12").Value = '''55.16'
$ws.Range("E12").Value = '  +13.95%  '
$ws.Range("E13").Value = '  -4.68%  '
$ws.Range("D14").Value = '''10.61'
$ws.Range("E14").Value = '  -5.18%  '
$ws.Range("D15").Value = '4.576.44'
$ws.Range("E15").Value = '  -2.50%  '
$ws.Range("D16").Value = '3.944.41'
$ws.Range("E16").Value = '  -3.26%  '
$ws.Range("D17").Value = '''13.91'
$ws.Range("E17").Value = '  -3.01%  '
$ws.Range("E18").Value = '  -3.80%  '
$ws.Range("E19").Value = '  -1.71%  '
$ws.Range("D21").Value = '70.928.85'
$ws.Range("E21").Value = '  -1.94%  '
$ws.Range("D22").Value = '''421.32'
$ws.Range("E22").Value = '  -5.68%  '
$ws.Range("E23").Value = '  -0.55%  '
$ws.Range("D24").Value = '''97.53'
$ws.Range("E24").Value = '  -7.38%  '
$ws.Range("E25").Value = '  +4.33%  '
$ws.Range("D26").Value = '''14.46'
$ws.Range("E26").Value = '  -5.25%  '
$ws.Range("D27").Value = '''11.39'
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("D28").Value = '''3.80'
$ws.Range("E28").Value = '  +15.37%  '
$ws.Range("D29").Value = '''10.69'
$ws.Range("E29").Value = '  -4.52%  '
$ws.Range("D30").Value = '''5.89'
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").Value = '''36.44'
$ws.Range("E31").Value = '  -4.50%  '
$ws.Range("D32").Value = '''7.90'
$ws.Range("E32").Value = '  +17.39%  '
$ws.Range("D33").Value = '''51.02'
$ws.Range("E33").Value = '  +19.11%  '
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("D35").Value = '''13.33'
$ws.Range("E35").Value = '  -3.53%  '
$ws.Range("D36").Value = '''682.93'
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("D37").Value = '''65.47'
$ws.Range("E37").Value = '  -3.63%  '
$ws.Range("D38").Value = '''0.441'
$ws.Range("E38").Value = '  +1.86%  '
$ws.Range("D39").Value = '0.0₃0815'
$ws.Range("E39").Value = '  -7.16%  '
$ws.Range("E40").Value = '  -3.48%  '
$ws.Range("E41").Value = '  -3.53%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("E44").Value = '  -4.02%  '
$ws.Range("D45").Value = '''3.18'
$ws.Range("E45").Value = '  -1.20%  '
$ws.Range("D46").Value = '''10.03'
$ws.Range("E46").Value = '  +2.49%  '
$ws.Range("E47").Value = '  -5.56%  '
$ws.Range("D48").Value = '''2.66'
$ws.Range("E48").Value = '  -2.46%  '
$ws.Range("E49").Value = '  -2.39%  '
$ws.Range("E50").Value = '  -2.56%  '
$ws.Range("D51").Value = '''144.23'
$ws.Range("E51").Value = '  -0.44%  '
